$d = $word.ActiveDocument

# Locate the current last paragraph of the document (the end of the
# trailing run of empty "PargrafodaLista" list paragraphs, just before
# the final section properties).
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Insert a new paragraph after it. Word (and this COM host) carries the
# style/numbering of the preceding paragraph onto the freshly inserted
# one automatically.
$insertionPoint = $lastPara.Range.Duplicate
$insertionPoint.Collapse(0)
$insertionPoint.InsertParagraphAfter()

# Re-fetch the brand-new (now last) paragraph and normalise its OOXML so
# it mirrors the other empty list paragraphs exactly: just a <w:pPr> with
# the "PargrafodaLista" style and numbering (ilvl 0 / numId 10), no
# stray empty run.
$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount)
$newRange = $newPara.Range.Duplicate
$null = $newRange.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='PargrafodaLista'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='10'/></w:numPr></w:pPr></w:p>")

$d.Save()
